$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "ODI Batting" sheet: remove the empty B2 cell (it currently is an
#    empty inline-string cell; clearing its value drops it entirely).
# ------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B2").Value = $null

# ------------------------------------------------------------------
# 2. Add a brand new worksheet "ODI Batting Extra" as the last sheet.
# ------------------------------------------------------------------
$lastIndex  = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($lastIndex)
$wsExtra    = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

# ------------------------------------------------------------------
# 3. Header row - reuse the same bold/centered/bordered style that is
#    already used for header rows on the other sheets.
# ------------------------------------------------------------------
$wsPlayerInfo = $wb.Worksheets.Item("Player Info")
$wsPlayerInfo.Range("A1:D1").Copy()
$wsExtra.Range("A1:F1").PasteSpecial(-4122)

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

# ------------------------------------------------------------------
# 4. Data rows. MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are
#    stored as text (matching how the rest of the workbook stores its
#    "numeric looking" values), BATTING_POSITION is a real number.
#    Rows 2 and 4 leave BATTING_POSITION/NUM_4/NUM_6/PERCENT blank.
# ------------------------------------------------------------------
$wsExtra.Range("A2:A4").NumberFormat = "@"
$wsExtra.Range("B2").NumberFormat = "@"
$wsExtra.Range("B4").NumberFormat = "@"
$wsExtra.Range("C2:D4").NumberFormat = "@"
$wsExtra.Range("E2:E4").NumberFormat = "@"

$wsExtra.Range("A2").Value = "4686"
$wsExtra.Range("B2").Value = ""
$wsExtra.Range("C2").Value = ""
$wsExtra.Range("D2").Value = ""
$wsExtra.Range("E2").Value = ""
$wsExtra.Range("F2").Value = "NO"

$wsExtra.Range("A3").Value = "4688"
$wsExtra.Range("B3").Value = 8
$wsExtra.Range("C3").Value = "1"
$wsExtra.Range("D3").Value = "1"
$wsExtra.Range("E3").Value = "6.59%"
$wsExtra.Range("F3").Value = "NO"

$wsExtra.Range("A4").Value = "4690"
$wsExtra.Range("B4").Value = ""
$wsExtra.Range("C4").Value = ""
$wsExtra.Range("D4").Value = ""
$wsExtra.Range("E4").Value = ""
$wsExtra.Range("F4").Value = "NO"
